# Minor changes to Light power
#
# 1) LightSourcePowerKeyMeasurements: remove the "power_set_point_pct" column
#    (was column B), shifting all following columns one place to the left.
# 2) PowerSample: rename/insert columns so the sheet becomes
#    acquisition_datetime, light_source, sampling_location, power_set_point,
#    power_mw, integration_time_ms - and add a list data validation (dropdown)
#    on the new "sampling_location" column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) LightSourcePowerKeyMeasurements
# ---------------------------------------------------------------------------
$wsKM = $wb.Worksheets.Item("LightSourcePowerKeyMeasurements")

# Column B was "power_set_point_pct" - delete the whole column, which shifts
# power_mean_mw, power_median_mw, power_std_mw, power_min_mw, power_max_mw,
# linearity, table_data, data_reference, linked_references, name, description
# one column to the left (C->B ... M->L).
$wsKM.Range("B1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2) PowerSample
# ---------------------------------------------------------------------------
$wsPS = $wb.Worksheets.Item("PowerSample")

# Before: A=light_source, B=sampling_datetime, C=power_mw

# Insert a new column at the front for "acquisition_datetime".
# This shifts light_source -> B, sampling_datetime -> C, power_mw -> D.
$wsPS.Range("A1").EntireColumn.Insert()
$wsPS.Cells.Item(1, 1).Value = "acquisition_datetime"

# Rename former "sampling_datetime" (now column C) to "sampling_location".
$wsPS.Cells.Item(1, 3).Value = "sampling_location"

# Insert a new column before power_mw (currently column D) for "power_set_point".
# This shifts power_mw -> E.
$wsPS.Range("D1").EntireColumn.Insert()
$wsPS.Cells.Item(1, 4).Value = "power_set_point"

# Add the trailing "integration_time_ms" column.
$wsPS.Cells.Item(1, 6).Value = "integration_time_ms"

# Add a dropdown list validation on the sampling_location column (C), for all
# data rows.
$rngLoc = $wsPS.Range("C2:C1048576")
$rngLoc.Validation.Add(3, 1, 3, """SOURCE_EXIT,FIBER_EXIT,OBJECTIVE_BACKFOCAL,OBJECTIVE_EXIT,OBJECTIVE_FOCAL,OTHER""")
$rngLoc.Validation.IgnoreBlank = 1
$rngLoc.Validation.InCellDropdown = 1
$rngLoc.Validation.ShowInput = 0
$rngLoc.Validation.ShowError = 0
